$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the KS4 / KS5 destinations labels (removed the "- provisional" suffix)
$ws.Range("A11").Value = "Key Stage 4 (KS4) destinations"
$ws.Range("A12").Value = "Key Stage 5 (KS5) destinations"

# Load in the new National Pupil Database permalinks for KS4 / KS5 destinations
$ws.Range("B11").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/data-tables/permalink/1683bef5-5daa-49d7-9323-08db08498a11'>National Pupil Database</a>"
$ws.Range("B12").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/data-tables/permalink/1703fe2c-2e57-4bfe-9325-08db08498a11'>National Pupil Database</a>"
